$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 132, shifting existing rows 132.. down by one.
$ws.Rows.Item(132).Insert()

# Populate new row 132 with data (copy of old row132 pattern, but with updated values)
$ws.Cells.Item(132, 1).Value = 10
$ws.Cells.Item(132, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(132, 3).Value = "La Araucanía"
$ws.Cells.Item(132, 4).Value = 44587
$ws.Cells.Item(132, 5).Value = 9
$ws.Cells.Item(132, 6).Value = 100112009
$ws.Cells.Item(132, 7).Value = "Acelga"
$ws.Cells.Item(132, 8).Value = "Sin especificar"
$ws.Cells.Item(132, 9).Value = "Primera"
$ws.Cells.Item(132, 10).Value = 55
$ws.Cells.Item(132, 11).Value = 7000
$ws.Cells.Item(132, 12).Value = 8000
$ws.Cells.Item(132, 13).Value = 7455
$ws.Cells.Item(132, 14).Value = "`$/docena de atados (12 kilos)"
$ws.Cells.Item(132, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(132, 16).Value = 621
$ws.Cells.Item(132, 17).Value = 12
$ws.Cells.Item(132, 18).Value = "Hortaliza"
